$wb = $excel.ActiveWorkbook

# ---------- About sheet ----------
$about = $wb.Worksheets.Item("About")

# Insert a new row at 10 (shifts old rows 10-17 down to 11-18).
# Formulas referencing About!$A$16 elsewhere auto-update to About!$A$17.
$about.Rows.Item(10).Insert()

# Row 1: remove the trailing empty B1 cell entirely
$about.Range("B1").Clear()

# Row 3: B3 text stays the same, but drop its explicit style (becomes default/Normal)
$about.Range("B3").Style = "Normal"

# Row 4: A4 was an empty styled cell -> clear it entirely
$about.Range("A4").Clear()

# Row 5: A5 empty styled cell -> clear; B5 text same, style -> Normal
$about.Range("A5").Clear()
$about.Range("B5").Style = "Normal"

# Row 6: A6 empty styled cell -> clear; B6 keeps its style (hyperlink-like style s=3)
$about.Range("A6").Clear()

# Row 7: A7 empty styled cell -> clear; B7 text same, style -> Normal
$about.Range("A7").Clear()
$about.Range("B7").Style = "Normal"

# Row 9: B9 was an empty styled cell -> clear it entirely
$about.Range("B9").Clear()

# Row 10 (newly inserted): new note text, default style
$about.Range("A10").Value2 = "When considering the Social Cost of Carbon, meant to capture the long-term economic damage caused by one"
$about.Range("A10").Style = "Normal"
$about.Range("B10").Clear()

# Row 11 (shifted from old row 10): split off second half of the note text, keep old style
$about.Range("A11").Value2 = "ton of carbon dioxide emitted, the U.S. government typically uses the figures based on"
$about.Range("B11").Clear()

# Row 12 (shifted from old row 11): "a 3% discount rate..." -> style becomes Normal
$about.Range("A12").Style = "Normal"
$about.Range("B12").Clear()

# Row 13 (shifted from old row 12): "Note that this differs..." -> style becomes Normal
$about.Range("A13").Style = "Normal"
$about.Range("B13").Clear()

# Row 14 (shifted from old row 13): "based on the Office of Management..." -> style becomes Normal
$about.Range("A14").Style = "Normal"
$about.Range("B14").Clear()

# Row 16 (shifted from old row 15): "We adjust 2007 dollars..." keeps its style; clear B16
$about.Range("B16").Clear()

# Row 17 (shifted from old row 16): 1.109 keeps its style; clear B17
$about.Range("B17").Clear()

# Row 18 (shifted from old row 17): "See cpi.xlsx..." keeps its style (only cell in this row now)

# ---------- SourceData sheet ----------
$src = $wb.Worksheets.Item("SourceData")
# Column F was a set of unused, empty styled placeholder cells -> remove entirely
$src.Range("F1:F44").Clear()

# ---------- SCoC sheet ----------
$scoc = $wb.Worksheets.Item("SCoC")
# Header label now reflects the gram-based unit
$scoc.Range("B1").Value2 = "Social Cost of Carbon ($/g CO2e)"
# Column A values/formulas unchanged, but drop the explicit (redundant) style
$scoc.Range("A2:A42").Style = "Normal"

Write-Host "done"
